# "Generate Report for Archive"
# - Update localization status text: "Ready for handoff" -> "In Translation"
#   (every cell that shows this status, across all three sheets)
# - Shrink the now-narrower "Status"/per-locale status columns to match the
#   shorter text: Overview!E:F and the "Status" column (col C) on the
#   zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$Overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# --- Overview sheet: columns E (zh-cn) and F (de-de), rows 2-3 ---
$Overview.Range("E2").Value = $newStatus
$Overview.Range("F2").Value = $newStatus
$Overview.Range("E3").Value = $newStatus
$Overview.Range("F3").Value = $newStatus

# --- zh-cn / de-de sheets: "Status" column (C), rows 2-3 ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# --- Re-fit the status columns now that the text is shorter ---
# (ColumnWidth is expressed in characters; 12.5 is the setting this host
# resolves to the tighter width matching the shorter "In Translation" text.)
$Overview.Columns.Item(5).ColumnWidth = 12.5
$Overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
